# Applies the "merging with new 2 test cases" edit described by the diff:
#  - TC01 (sheet1): replace the actiTIME login-sample row with a Zoho CRM
#    login-sample row, turn A2 into a mailto hyperlink, widen columns A & C.
#  - TC03 (sheet3): replace the Version/actiTIME build row with a
#    campaign-title/Zoho CRM row, widen column A, make TC03 the active tab.
#  - TC02 / TC04 keep their textual content (only shared-string ids shift,
#    which happens automatically as strings are added/removed elsewhere).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# TC01
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TC01")

# New login-sample content (was admin / manager / actiTIME - Login / actiTIME - Enter Time-Track)
$ws1.Range("A2").Value = "rashmi@dell.com"
$ws1.Range("B2").Value = 123456
$ws1.Range("C2").Value = "Zoho CRM - Sign in"
$ws1.Range("D2").Value = "Zoho CRM - Home Page"

# Turn the user name cell into a mailto hyperlink (adopts the built-in
# "Hyperlink" cell style automatically).
[void]$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:rashmi@dell.com")

# Selection moves from A29 to D2.
$ws1.Range("D2").Select()

# Column widths grow to fit the new, longer content.
$ws1.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws1.Columns.Item(3).ColumnWidth = 16.666666666666668

# ---------------------------------------------------------------------
# TC03
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("TC03")

# New campaign-title sample content (was Version / actiTIME 2019.Pro).
$ws3.Range("A1").Value = "campaign title"
$ws3.Range("A2").Value = "Zoho CRM - Create Campaign"

$ws3.Columns.Item(1).ColumnWidth = 26.166666666666668

# Selection moves from A2 to E12.
$ws3.Range("E12").Select()

# TC03 becomes the active / selected sheet (activeTab=2, tabSelected on TC03).
$ws3.Activate()
